$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.644.93'
$ws.Range("E2").Value = '  +0.64%  '

# Row 3
$ws.Range("D3").Value = '1.847.22'
$ws.Range("E3").Value = '  -0.14%  '

# Row 4
$ws.Range("D4").Value = '''1.035'
$ws.Range("E4").Value = '  +0.29%  '

# Row 5
$ws.Range("D5").Value = '''322.28'
$ws.Range("E5").Value = '  +0.34%  '

# Row 6
$ws.Range("D6").Value = '''1.032'
$ws.Range("E6").Value = '  +0.28%  '

# Row 7
$ws.Range("D7").Value = '''0.4373'
$ws.Range("E7").Value = '  +0.13%  '

# Row 8
$ws.Range("D8").Value = '''0.3785'
$ws.Range("E8").Value = '  +0.95%  '

# Row 9
$ws.Range("D9").Value = '''0.07350'
$ws.Range("E9").Value = '  -0.46%  '

# Row 10
$ws.Range("D10").Value = '''0.8766'
$ws.Range("E10").Value = '  +0.33%  '

# Row 11
$ws.Range("D11").Value = '''21.31'
$ws.Range("E11").Value = '  -0.04%  '

# Row 12
$ws.Range("D12").Value = '1.851.23'
$ws.Range("E12").Value = '  -0.91%  '

# Row 13
$ws.Range("D13").Value = '''5.464'
$ws.Range("E13").Value = '  -0.82%  '

# Row 14
$ws.Range("D14").Value = '''6.684'
$ws.Range("E14").Value = '  +0.28%  '

# Row 15
$ws.Range("D15").Value = '''0.07197'
$ws.Range("E15").Value = '  -0.20%  '

# Row 16
$ws.Range("D16").Value = '''85.42'
$ws.Range("E16").Value = '  +3.51%  '

# Row 17
$ws.Range("D17").Value = '''1.038'
$ws.Range("E17").Value = '  +0.38%  '

# Row 18
$ws.Range("D18").Value = '''0.000009018'
$ws.Range("E18").Value = '  +0.17%  '

# Row 19
$ws.Range("D19").Value = '''1.032'
$ws.Range("E19").Value = '  +0.33%  '

# Row 20
$ws.Range("D20").Value = '''15.40'
$ws.Range("E20").Value = '  -0.01%  '

# Row 21
$ws.Range("D21").Value = '27.690.57'
$ws.Range("E21").Value = '  +0.68%  '

# Row 22
$ws.Range("D22").Value = '''5.257'
$ws.Range("E22").Value = '  +0.23%  '

# Row 23
$ws.Range("D23").Value = '''11.13'
$ws.Range("E23").Value = '  -0.57%  '

# Row 24
$ws.Range("D24").Value = '2.083.73'
$ws.Range("E24").Value = '  +0.24%  '

# Row 25
$ws.Range("D25").Value = '''2.068'
$ws.Range("E25").Value = '  +7.24%  '

# Row 26
$ws.Range("D26").Value = '''158.77'
$ws.Range("E26").Value = '  +0.78%  '

# Row 27
$ws.Range("D27").Value = '''18.58'
$ws.Range("E27").Value = '  -0.51%  '

# Row 28
$ws.Range("D28").Value = '''5.303'
$ws.Range("E28").Value = '  +0.77%  '

# Row 29
$ws.Range("D29").Value = '''1.971'
$ws.Range("E29").Value = '  +1.82%  '

# Row 30
$ws.Range("D30").Value = '''118.50'
$ws.Range("E30").Value = '  +1.45%  '

# Row 31
$ws.Range("D31").Value = '''0.09088'
$ws.Range("E31").Value = '  +0.59%  '

# Row 32
$ws.Range("B32").Value = 'HuobiToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D32").Value = '''3.050'
$ws.Range("E32").Value = '  +5.73%  '

# Row 33
$ws.Range("D33").Value = '''1.196'
$ws.Range("E33").Value = '  -0.46%  '

# Row 34
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = '''0.7573'
$ws.Range("E34").Value = '  -0.50%  '

# Row 35
$ws.Range("D35").Value = '''4.532'
$ws.Range("E35").Value = '  +0.78%  '

# Row 36
$ws.Range("D36").Value = '''1.034'
$ws.Range("E36").Value = '  +0.40%  '

# Row 37
$ws.Range("D37").Value = '''1.147'
$ws.Range("E37").Value = '  -0.01%  '

# Row 38
$ws.Range("D38").Value = '''0.01971'
$ws.Range("E38").Value = '  -0.05%  '

# Row 39
$ws.Range("D39").Value = '''0.05260'
$ws.Range("E39").Value = '  -0.34%  '

# Row 40
$ws.Range("D40").Value = '''2.836'
$ws.Range("E40").Value = '  +1.09%  '

# Row 41
$ws.Range("D41").Value = '''0.5161'
$ws.Range("E41").Value = '  +0.19%  '

# Row 42
$ws.Range("D42").Value = '''6.949'
$ws.Range("E42").Value = '  +3.69%  '

# Row 43
$ws.Range("D43").Value = '''0.1662'
$ws.Range("E43").Value = '  -0.26%  '

# Row 44
$ws.Range("D44").Value = '''8.635'
$ws.Range("E44").Value = '  +1.21%  '

# Row 45
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value = '''110.14'
$ws.Range("E45").Value = '  +1.21%  '

# Row 46
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '''10.68'
$ws.Range("E46").Value = '  +1.54%  '

# Row 47
$ws.Range("D47").Value = '''1.035'
$ws.Range("E47").Value = '  +0.46%  '

# Row 48
$ws.Range("D48").Value = '''0.06522'
$ws.Range("E48").Value = '  +2.09%  '

# Row 49
$ws.Range("D49").Value = '''1.700'
$ws.Range("E49").Value = '  -0.53%  '

# Row 50
$ws.Range("D50").Value = '''0.4680'
$ws.Range("E50").Value = '  +0.87%  '

# Row 51
$ws.Range("D51").Value = '''1.860'
$ws.Range("E51").Value = '  -0.14%  '
